$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NroPoliza (E2) - change to new text value (forced text, same visual style)
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "12112002379"
$ws.Range("E2").NumberFormat = "General"

# Update FechaSiniestro (G2) - change date text (apostrophe keeps quote-prefixed text style)
$ws.Range("G2").Value = "'27/05/2021"

# Update selection
$ws.Range("G3").Select()
